# Adiciona arquivos com informações do plot do PCR com filtro SG
# Reorders header columns (R², RMSE, Offset, Slope), reorders/updates data
# columns to match (C=R², D=RMSE, E=Offset, F=Slope), fixes several
# Attribute/Y mislabels, updates many numeric values, and appends new
# "Validação" rows for each attribute.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Attribute"
$ws.Range("B1").Value = "Y"
$ws.Range("C1").Value = "R²"
$ws.Range("D1").Value = "RMSE"
$ws.Range("E1").Value = "Offset"
$ws.Range("F1").Value = "Slope"

# --- Data rows (Attribute, Y, R², RMSE, Offset, Slope) ---
$data = @(
    @("SST",         "Referência", 0.7652488218511146, 1.287430089288877, 3.287053068205883, 0.7652488218511148),
    @("SST",         "Predição",   0.7352005490747384, 1.367345623484522, 3.482316337791213, 0.7521081881774457),
    @("SST",         "Validação",  0.7348733649008815, 1.08223076245959,  1.832791360085006,  0.8698628171401067),
    @("PH",          "Referência", 0.4688054746938572, 0.2251811710418577,1.743119387888039,  0.4688054746938569),
    @("PH",          "Predição",   0.390448746200474,  0.2412183778558815,1.841827598863998,  0.438344113288157),
    @("PH",          "Validação",  0.5551700929815269, 0.174563123106224, 0.5250335691233268, 0.8500194321630428),
    @("AT",          "Referência", 0.5529044951640683, 0.3922437331892594,0.5015287438418431, 0.5529044951640681),
    @("AT",          "Predição",   0.5022059794686908, 0.4138859651465448,0.5249115170629447, 0.529348735615595),
    @("AT",          "Validação",  0.5884635316029205, 0.2802708219034881,0.05050495057193216,0.9210674899942144),
    @("FIRMEZA (N)", "Referência", 0.5497937355681533, 68.95711596316038, 230.9589342261019,  0.5497937355681534),
    @("FIRMEZA (N)", "Predição",   0.4730347656370216, 74.60437642718099, 244.9523546846835,  0.5215083906750153),
    @("FIRMEZA (N)", "Validação",  0.4133020347533339, 59.17165279768364, 102.7664707452828,  0.7823503012261875),
    @("UBS (%)",     "Referência", 0.6704659607954542, 1.846388222751209, 5.057729332751462,  0.6704659607954542),
    @("UBS (%)",     "Predição",   0.6334989672813065, 1.947199645618298, 5.286829133664152,  0.6562443076208077),
    @("UBS (%)",     "Validação",  0.7723798644787406, 1.033875521911286, 2.705202574476647,  0.8201895022104482)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
